# Update countries & provincias Spain
# - Swap the display order of "Nueva Zelanda" / "Tunez" (row 137 / 138)
# - Refresh the "Datos actualizados" timestamp
# - Refresh the COVID-19 stat counters for a number of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names in rows 137 / 138 (Tunez now listed before Nueva Zelanda) ---
$ws.Cells.Item(137, 1).Value = "Tunez"
$ws.Cells.Item(138, 1).Value = "Nueva Zelanda"

# --- Update "last updated" banner text ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 4 de Agosto de 2020 a las 20:32"

# --- Refresh numeric stats, row by row (columns B..H = totals/new/active/recovered/critical/deaths-today/deaths) ---

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 4887648
$ws.Cells.Item(4, 3).Value = 23732
$ws.Cells.Item(4, 4).Value = 2455199
$ws.Cells.Item(4, 5).Value = 2272853
$ws.Cells.Item(4, 7).Value = 668
$ws.Cells.Item(4, 8).Value = 159596

# Row 6 - India
$ws.Cells.Item(6, 2).Value = 1906520
$ws.Cells.Item(6, 3).Value = 51189
$ws.Cells.Item(6, 4).Value = 1280667
$ws.Cells.Item(6, 5).Value = 586033
$ws.Cells.Item(6, 7).Value = 849
$ws.Cells.Item(6, 8).Value = 39820

# Row 15 - Reino Unido
$ws.Cells.Item(15, 7).Value = 89
$ws.Cells.Item(15, 8).Value = 46299

# Row 20 - Turquia
$ws.Cells.Item(20, 2).Value = 234934
$ws.Cells.Item(20, 3).Value = 1083
$ws.Cells.Item(20, 4).Value = 218491
$ws.Cells.Item(20, 5).Value = 10678
$ws.Cells.Item(20, 7).Value = 18
$ws.Cells.Item(20, 8).Value = 5765

# Row 23 - Francia
$ws.Cells.Item(23, 2).Value = 192334
$ws.Cells.Item(23, 3).Value = 1039
$ws.Cells.Item(23, 5).Value = 79874

# Row 31 - Ecuador
$ws.Cells.Item(31, 4).Value = 70985
$ws.Cells.Item(31, 5).Value = 11170

# Row 60 - Argelia
$ws.Cells.Item(60, 2).Value = 32504
$ws.Cells.Item(60, 3).Value = 532
$ws.Cells.Item(60, 4).Value = 22375
$ws.Cells.Item(60, 5).Value = 8881
$ws.Cells.Item(60, 7).Value = 9
$ws.Cells.Item(60, 8).Value = 1248

# Row 65 - Moldavia
$ws.Cells.Item(65, 4).Value = 18167
$ws.Cells.Item(65, 5).Value = 6837

# Row 107 - Maldivas
$ws.Cells.Item(107, 2).Value = 4446
$ws.Cells.Item(107, 3).Value = 153
$ws.Cells.Item(107, 4).Value = 2693
$ws.Cells.Item(107, 5).Value = 1734
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 19

# Row 108 - Malaui
$ws.Cells.Item(108, 2).Value = 4273
$ws.Cells.Item(108, 3).Value = 1
$ws.Cells.Item(108, 4).Value = 2020
$ws.Cells.Item(108, 5).Value = 2130

# Row 117 - Mayotte
$ws.Cells.Item(117, 2).Value = 3023
$ws.Cells.Item(117, 3).Value = 15
$ws.Cells.Item(117, 5).Value = 246

# Row 118 - Suazilandia
$ws.Cells.Item(118, 2).Value = 2856
$ws.Cells.Item(118, 3).Value = 18
$ws.Cells.Item(118, 4).Value = 1258
$ws.Cells.Item(118, 5).Value = 1549
$ws.Cells.Item(118, 7).Value = 4
$ws.Cells.Item(118, 8).Value = 49

# Row 136 - Yemen
$ws.Cells.Item(136, 2).Value = 1760
$ws.Cells.Item(136, 3).Value = 26
$ws.Cells.Item(136, 5).Value = 391
$ws.Cells.Item(136, 7).Value = 7
$ws.Cells.Item(136, 8).Value = 506

# Row 137 - now Tunez (after the name swap above)
$ws.Cells.Item(137, 2).Value = 1584
$ws.Cells.Item(137, 3).Value = 19
$ws.Cells.Item(137, 4).Value = 1227
$ws.Cells.Item(137, 5).Value = 306
$ws.Cells.Item(137, 8).Value = 51

# Row 138 - now Nueva Zelanda (after the name swap above)
$ws.Cells.Item(138, 2).Value = 1567
$ws.Cells.Item(138, 4).Value = 1523
$ws.Cells.Item(138, 5).Value = 22
$ws.Cells.Item(138, 8).Value = 22

# Row 151 - Republica del Chad
$ws.Cells.Item(151, 2).Value = 938
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 5).Value = 49

# Row 157 - Lesoto
$ws.Cells.Item(157, 2).Value = 726
$ws.Cells.Item(157, 3).Value = 8
$ws.Cells.Item(157, 4).Value = 174
$ws.Cells.Item(157, 5).Value = 531
$ws.Cells.Item(157, 7).Value = 2
$ws.Cells.Item(157, 8).Value = 21

# Row 186 - Monaco
$ws.Cells.Item(186, 2).Value = 123
$ws.Cells.Item(186, 3).Value = 2
$ws.Cells.Item(186, 5).Value = 14
